# Insert a new "panel" row into the example1 sheet's properties table,
# right after the "nSpans" row (row 2), pushing all subsequent rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("example1")

# Insert a new row before current row 3, shifting existing rows 3.. down.
$ws.Rows.Item(3).Insert()

# Populate the new row 3 with the panel property.
$ws.Cells.Item(3, 1).Value = "panel"
$ws.Cells.Item(3, 2).Value = "int"
$ws.Cells.Item(3, 3).Value = "string"
$ws.Cells.Item(3, 4).Value = "Panel location of section [interior/end]"

$ws.Range("J13").Select()
